# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.480.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5358"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2924"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06954"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.850.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7242"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07195"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.971"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007870"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.490.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.080.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.581"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.974"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.165"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.708"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.151"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.241"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08881"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.023"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04830"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7236"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.130"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.094"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.294"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4652"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9027"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.866"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.401"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.025"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1241"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4031"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8909"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05741"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "
